$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price/volume snapshot (Coin, Link, Price, Volume(1h)).
# For Price (column D) values that look like plain numbers (e.g. "333.42"),
# a leading apostrophe is used so Excel keeps them as text (quote-prefixed),
# matching the original inline-string text values rather than converting them
# to numeric cells.

# Row 2
$ws.Range("D2").Value = '29.034.80'
$ws.Range("E2").Value = '  +2.36%  '

# Row 3
$ws.Range("D3").Value = '1.906.36'
$ws.Range("E3").Value = '  +2.01%  '

# Row 4
$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = '  -0.14%  '

# Row 5
$ws.Range("D5").Value = "'333.42"
$ws.Range("E5").Value = '  -1.66%  '

# Row 6
$ws.Range("D6").Value = "'0.9997"
$ws.Range("E6").Value = '  -0.13%  '

# Row 7
$ws.Range("D7").Value = "'0.4634"
$ws.Range("E7").Value = '  -1.49%  '

# Row 8
$ws.Range("D8").Value = "'0.4084"
$ws.Range("E8").Value = '  +2.93%  '

# Row 9
$ws.Range("D9").Value = "'47.92"
$ws.Range("E9").Value = '  +0.82%  '

# Row 10
$ws.Range("D10").Value = "'0.08015"
$ws.Range("E10").Value = '  -0.34%  '

# Row 11
$ws.Range("D11").Value = "'1.006"
$ws.Range("E11").Value = '  +0.48%  '

# Row 12
$ws.Range("D12").Value = "'21.78"
$ws.Range("E12").Value = '  -0.70%  '

# Row 13
$ws.Range("D13").Value = '1.907.03'
$ws.Range("E13").Value = '  +1.58%  '

# Row 14
$ws.Range("D14").Value = "'5.948"
$ws.Range("E14").Value = '  -1.54%  '

# Row 15
$ws.Range("D15").Value = "'7.096"
$ws.Range("E15").Value = '  -2.17%  '

# Row 16
$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = '  -0.24%  '

# Row 17
$ws.Range("D17").Value = "'88.96"
$ws.Range("E17").Value = '  -2.48%  '

# Row 18
$ws.Range("D18").Value = "'0.00001032"
$ws.Range("E18").Value = '  -1.07%  '

# Row 19
$ws.Range("D19").Value = "'0.06569"
$ws.Range("E19").Value = '  -0.80%  '

# Row 20
$ws.Range("D20").Value = "'17.54"
$ws.Range("E20").Value = '  -0.10%  '

# Row 21
$ws.Range("D21").Value = "'0.9998"
$ws.Range("E21").Value = '  -0.09%  '

# Row 22
$ws.Range("D22").Value = '29.034.03'
$ws.Range("E22").Value = '  +2.31%  '

# Row 23
$ws.Range("D23").Value = "'5.450"
$ws.Range("E23").Value = '  -0.53%  '

# Row 24
$ws.Range("E24").Value = '  +1.94%  '

# Row 25
$ws.Range("E25").Value = '  -0.78%  '

# Row 26
$ws.Range("D26").Value = '2.133.90'
$ws.Range("E26").Value = '  +1.91%  '

# Row 27
$ws.Range("D27").Value = "'157.75"
$ws.Range("E27").Value = '  -1.92%  '

# Row 28
$ws.Range("D28").Value = "'19.73"
$ws.Range("E28").Value = '  -0.12%  '

# Row 29
$ws.Range("D29").Value = "'2.102"
$ws.Range("E29").Value = '  -0.95%  '

# Row 30
$ws.Range("D30").Value = "'5.408"

# Row 31
$ws.Range("D31").Value = "'118.95"
$ws.Range("E31").Value = '  -1.19%  '

# Row 32
$ws.Range("D32").Value = "'0.9810"
$ws.Range("E32").Value = '  +1.07%  '

# Row 33
$ws.Range("D33").Value = "'0.09424"
$ws.Range("E33").Value = '  -0.92%  '

# Row 34
$ws.Range("D34").Value = "'1.422"
$ws.Range("E34").Value = '  +3.55%  '

# Row 35
$ws.Range("D35").Value = "'3.589"
$ws.Range("E35").Value = '  -0.19%  '

# Row 36
$ws.Range("D36").Value = "'5.317"
$ws.Range("E36").Value = '  -0.61%  '

# Row 37
$ws.Range("D37").Value = "'0.06092"
$ws.Range("E37").Value = '  -0.14%  '

# Row 38
$ws.Range("D38").Value = "'0.02242"
$ws.Range("E38").Value = '  -0.56%  '

# Row 39
$ws.Range("D39").Value = "'8.384"
$ws.Range("E39").Value = '  -0.02%  '

# Row 40
$ws.Range("D40").Value = "'1.175"
$ws.Range("E40").Value = '  -1.04%  '

# Row 41
$ws.Range("D41").Value = "'0.5816"
$ws.Range("E41").Value = '  -2.37%  '

# Row 42
$ws.Range("D42").Value = "'0.9995"
$ws.Range("E42").Value = '  -0.07%  '

# Row 43
$ws.Range("D43").Value = "'10.19"
$ws.Range("E43").Value = '  -1.77%  '

# Row 44
$ws.Range("D44").Value = "'0.1824"
$ws.Range("E44").Value = '  -2.73%  '

# Row 45
$ws.Range("E45").Value = '  -1.91%  '

# Row 46
$ws.Range("D46").Value = "'2.311"
$ws.Range("E46").Value = '  +11.01%  '

# Row 47
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = "'0.5511"
$ws.Range("E47").Value = '  -1.50%  '

# Row 48
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = "'12.07"
$ws.Range("E48").Value = '  -1.22%  '

# Row 49
$ws.Range("D49").Value = "'1.914"
$ws.Range("E49").Value = '  -2.24%  '

# Row 50
$ws.Range("D50").Value = "'48.49"
$ws.Range("E50").Value = '  +24.97%  '

# Row 51
$ws.Range("D51").Value = "'0.07030"
$ws.Range("E51").Value = '  +2.13%  '
